$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3331.7917
$ws.Range("I76").Value = 3181.5
$ws.Range("K76").Value = 3181.5
$ws.Range("M76").Value = -2866.5
$ws.Range("H79").Value = 3331.7917
$ws.Range("I79").Value = 3181.5
$ws.Range("K79").Value = 3181.5
$ws.Range("M79").Value = -2089.5
$ws.Range("H105").Value = 31738.584
$ws.Range("J105").Value = 31738.584
$ws.Range("L105").Value = 31738.584
$ws.Range("N105").Value = -38726.584
$ws.Range("H112").Value = 1303.75
$ws.Range("J112").Value = 1318.3636
$ws.Range("L112").Value = 3955.0908
$ws.Range("N112").Value = -6171.0908
$ws.Range("H129").Value = 836.37
$ws.Range("J129").Value = 864.6316
$ws.Range("L129").Value = 2593.8948
$ws.Range("N129").Value = -12593.8948
$ws.Range("H137").Value = 2167119.2
$ws.Range("I137").Value = 2382581.2
$ws.Range("K137").Value = 7147743.600000001
$ws.Range("M137").Value = -7145193.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4759.0933
$ws.Range("I32").Value = 3122.7693
$ws.Range("J32").Value = 9823.904
$ws.Range("K32").Value = 3122.7693
$ws.Range("L32").Value = 9823.904
$ws.Range("M32").Value = -2835.7693
$ws.Range("N32").Value = -10397.904
$ws.Range("H132").Value = 2409.0476
$ws.Range("I132").Value = 1026.6666
$ws.Range("J132").Value = 4252.222
$ws.Range("K132").Value = 3079.9998
$ws.Range("L132").Value = 12756.666
$ws.Range("M132").Value = -549.9998000000001
$ws.Range("N132").Value = -17816.666
$ws.Range("H137").Value = 41735
$ws.Range("J137").Value = 41735
$ws.Range("L137").Value = 41735
$ws.Range("N137").Value = -51935
$ws.Range("H139").Value = 41515.668
$ws.Range("J139").Value = 41515.668
$ws.Range("L139").Value = 41515.668
$ws.Range("N139").Value = -51795.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 40820.383
$ws.Range("J138").Value = 40820.383
$ws.Range("L138").Value = 40820.383
$ws.Range("N138").Value = -51100.383

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6604.5
$ws.Range("J31").Value = 12735.429
$ws.Range("L31").Value = 12735.429
$ws.Range("N31").Value = -13325.429
$ws.Range("H34").Value = 6604.5
$ws.Range("J34").Value = 12735.429
$ws.Range("L34").Value = 12735.429
$ws.Range("N34").Value = -13139.429
$ws.Range("H99").Value = 7695825.5
$ws.Range("I99").Value = 12501882
$ws.Range("J99").Value = 6135
$ws.Range("K99").Value = 12501882
$ws.Range("L99").Value = 6135
$ws.Range("M99").Value = -12500384
$ws.Range("N99").Value = -9131
$ws.Range("H105").Value = 2626.25
$ws.Range("J105").Value = 2666.6667
$ws.Range("L105").Value = 2666.6667
$ws.Range("N105").Value = -6160.6667
$ws.Range("H126").Value = 7695825.5
$ws.Range("I126").Value = 12501882
$ws.Range("J126").Value = 6135
$ws.Range("K126").Value = 37505646
$ws.Range("L126").Value = 18405
$ws.Range("M126").Value = -37503176
$ws.Range("N126").Value = -23345

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 581582.25
$ws.Range("I5").Value = 412.42856
$ws.Range("K5").Value = 1237.28568
$ws.Range("M5").Value = -1125.28568
$ws.Range("H80").Value = 8000
$ws.Range("I80").Value = 7000.6665
$ws.Range("J80").Value = 8749.5
$ws.Range("K80").Value = 21001.9995
$ws.Range("L80").Value = 26248.5
$ws.Range("M80").Value = -20065.9995
$ws.Range("N80").Value = -28120.5
$ws.Range("H83").Value = 8000
$ws.Range("I83").Value = 7000.6665
$ws.Range("J83").Value = 8749.5
$ws.Range("K83").Value = 63005.9985
$ws.Range("L83").Value = 78745.5
$ws.Range("M83").Value = -58325.9985
$ws.Range("N83").Value = -88105.5
$ws.Range("H92").Value = 710
$ws.Range("I92").Value = 455.55554
$ws.Range("K92").Value = 1366.66662
$ws.Range("M92").Value = -118.66662
$ws.Range("H126").Value = 2800
$ws.Range("I126").Value = 2800
$ws.Range("K126").Value = 8400
$ws.Range("M126").Value = -3460
$ws.Range("H131").Value = 791.3535000000001
$ws.Range("J131").Value = 832.4066
$ws.Range("L131").Value = 2497.2198
$ws.Range("N131").Value = -12577.2198
$ws.Range("H135").Value = 581582.25
$ws.Range("I135").Value = 412.42856
$ws.Range("K135").Value = 3711.85704
$ws.Range("M135").Value = -1176.85704
$ws.Range("H138").Value = 2211.6667
$ws.Range("I138").Value = 923.3333
$ws.Range("K138").Value = 2769.9999
$ws.Range("M138").Value = 2370.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 23878
$ws.Range("J46").Value = 23878
$ws.Range("L46").Value = 23878
$ws.Range("N46").Value = -24190
$ws.Range("H107").Value = 777.1111
$ws.Range("I107").Value = 589.1111
$ws.Range("J107").Value = 965.1111
$ws.Range("K107").Value = 589.1111
$ws.Range("L107").Value = 965.1111
$ws.Range("M107").Value = 1330.8889
$ws.Range("N107").Value = -4805.1111
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 40217.2
$ws.Range("J137").Value = 40217.2
$ws.Range("L137").Value = 40217.2
$ws.Range("N137").Value = -50417.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 7665.6665
$ws.Range("I23").Value = 5000
$ws.Range("J23").Value = 8998.5
$ws.Range("K23").Value = 5000
$ws.Range("L23").Value = 8998.5
$ws.Range("M23").Value = -4770
$ws.Range("N23").Value = -9458.5
$ws.Range("H40").Value = 4547.2
$ws.Range("I40").Value = 3800.5715
$ws.Range("K40").Value = 3800.5715
$ws.Range("M40").Value = -3664.5715
$ws.Range("H122").Value = 5689.5386
$ws.Range("I122").Value = 3996.2856
$ws.Range("J122").Value = 7665
$ws.Range("K122").Value = 11988.8568
$ws.Range("L122").Value = 22995
$ws.Range("M122").Value = -9538.856800000001
$ws.Range("N122").Value = -27895
$ws.Range("H132").Value = 5413.3057
$ws.Range("I132").Value = 3321.2144
$ws.Range("K132").Value = 9963.643199999999
$ws.Range("M132").Value = -7433.643199999999
$ws.Range("H133").Value = 26098.691
$ws.Range("J133").Value = 35383
$ws.Range("L133").Value = 35383
$ws.Range("N133").Value = -40443
$ws.Range("H139").Value = 44438.332
$ws.Range("J139").Value = 44438.332
$ws.Range("L139").Value = 44438.332
$ws.Range("N139").Value = -54718.332
$ws.Range("H140").Value = 61836.11
$ws.Range("J140").Value = 61836.11
$ws.Range("L140").Value = 61836.11
$ws.Range("N140").Value = -72196.11
$ws.Range("H141").Value = 31771
$ws.Range("J141").Value = 31771
$ws.Range("L141").Value = 31771
$ws.Range("N141").Value = -42131

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 3628.2856
$ws.Range("I4").Value = 799.6667
$ws.Range("K4").Value = 799.6667
$ws.Range("M4").Value = -686.6667
$ws.Range("H62").Value = 26000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31248
$ws.Range("H65").Value = 26000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156240
$ws.Range("H81").Value = 2718.889
$ws.Range("I81").Value = 2315.75
$ws.Range("J81").Value = 5944
$ws.Range("K81").Value = 4631.5
$ws.Range("L81").Value = 11888
$ws.Range("M81").Value = -3570.5
$ws.Range("N81").Value = -14010
$ws.Range("H84").Value = 2718.889
$ws.Range("I84").Value = 2315.75
$ws.Range("J84").Value = 5944
$ws.Range("K84").Value = 23157.5
$ws.Range("L84").Value = 59440
$ws.Range("M84").Value = -17853.5
$ws.Range("N84").Value = -70048
$ws.Range("H107").Value = 879.2
$ws.Range("J107").Value = 898
$ws.Range("L107").Value = 2694
$ws.Range("N107").Value = -6534
$ws.Range("H122").Value = 10528.333
$ws.Range("I122").Value = 9160
$ws.Range("J122").Value = 12238.75
$ws.Range("K122").Value = 27480
$ws.Range("L122").Value = 36716.25
$ws.Range("M122").Value = -25030
$ws.Range("N122").Value = -41616.25
$ws.Range("H126").Value = 2353.8667
$ws.Range("I126").Value = 1553.0714
$ws.Range("J126").Value = 3054.5625
$ws.Range("K126").Value = 4659.2142
$ws.Range("L126").Value = 9163.6875
$ws.Range("M126").Value = -2189.2142
$ws.Range("N126").Value = -14103.6875
$ws.Range("H136").Value = 4484.2964
$ws.Range("I136").Value = 3841.077
$ws.Range("K136").Value = 11523.231
$ws.Range("M136").Value = -8973.231
$ws.Range("H138").Value = 45332.668
$ws.Range("J138").Value = 45332.668
$ws.Range("L138").Value = 45332.668
$ws.Range("N138").Value = -55612.668
$ws.Range("H139").Value = 40854.5
$ws.Range("J139").Value = 40854.5
$ws.Range("L139").Value = 40854.5
$ws.Range("N139").Value = -51134.5
$ws.Range("H140").Value = 31033.572
$ws.Range("J140").Value = 31033.572
$ws.Range("L140").Value = 31033.572
$ws.Range("N140").Value = -41393.572
$ws.Range("H141").Value = 31857
$ws.Range("J141").Value = 31857
$ws.Range("L141").Value = 31857
$ws.Range("N141").Value = -42217
